# Regenerate orders with updated distance/size codes.
# Distances: D51 -> D55, D64 -> D69, D80 -> D86
# Sizes:     S30 -> S31  (S20, S25 unchanged)
# These codes appear as substrings inside many strings (Condition,
# Filename_Left, Filename_Right, Distance, Size columns), so use a
# substring ("contains") find/replace over the whole used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D51", "D55", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("D64", "D69", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("D80", "D86", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$used.Replace("S30", "S31", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
